$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 426.25
$ws.Range("I11").Value = 426.25
$ws.Range("K11").Value = 426.25
$ws.Range("M11").Value = -286.25
$ws.Range("H17").Value = 1589700.4
$ws.Range("I17").Value = 1013
$ws.Range("K17").Value = 3039
$ws.Range("M17").Value = -2871
$ws.Range("H33").Value = 19347.273
$ws.Range("I33").Value = 23500.223
$ws.Range("J33").Value = 659
$ws.Range("K33").Value = 23500.223
$ws.Range("L33").Value = 659
$ws.Range("M33").Value = -23271.223
$ws.Range("N33").Value = -1117
$ws.Range("H40").Value = 2002919.4
$ws.Range("I40").Value = 5001499
$ws.Range("J40").Value = 3866.3333
$ws.Range("K40").Value = 5001499
$ws.Range("L40").Value = 3866.3333
$ws.Range("M40").Value = -5001324
$ws.Range("N40").Value = -4216.3333
$ws.Range("H41").Value = 78.42856999999999
$ws.Range("I41").Value = 63.333332
$ws.Range("J41").Value = 89.75
$ws.Range("K41").Value = 63.333332
$ws.Range("L41").Value = 89.75
$ws.Range("M41").Value = 376.666668
$ws.Range("N41").Value = -969.75
$ws.Range("H74").Value = 4091.6667
$ws.Range("I74").Value = 2421.875
$ws.Range("K74").Value = 2421.875
$ws.Range("M74").Value = -1485.875
$ws.Range("H77").Value = 4091.6667
$ws.Range("I77").Value = 2421.875
$ws.Range("K77").Value = 12109.375
$ws.Range("M77").Value = -7429.375
$ws.Range("H92").Value = 353.6875
$ws.Range("I92").Value = 334.83334
$ws.Range("K92").Value = 334.83334
$ws.Range("M92").Value = 913.16666
$ws.Range("H97").Value = 1249.3334
$ws.Range("J97").Value = 1249.3334
$ws.Range("L97").Value = 3748.0002
$ws.Range("N97").Value = -4740.0002
$ws.Range("H112").Value = 6412258.5
$ws.Range("J112").Value = 6495508.5
$ws.Range("L112").Value = 19486525.5
$ws.Range("N112").Value = -19488741.5
$ws.Range("H129").Value = 41730.2
$ws.Range("H132").Value = 1902.25
$ws.Range("I132").Value = 1786.7368
$ws.Range("K132").Value = 5360.2104
$ws.Range("M132").Value = -2830.2104
$ws.Range("H137").Value = 3054.0667
$ws.Range("I137").Value = 3863.5557
$ws.Range("K137").Value = 11590.6671
$ws.Range("M137").Value = -9040.667099999999
$ws.Range("H138").Value = 6176504.5
$ws.Range("I138").Value = 666.4167
$ws.Range("K138").Value = 1999.2501
$ws.Range("M138").Value = 3140.7499

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 109931.664
$ws.Range("I31").Value = 17603.8
$ws.Range("J31").Value = 571571
$ws.Range("K31").Value = 17603.8
$ws.Range("L31").Value = 571571
$ws.Range("M31").Value = -17309.8
$ws.Range("N31").Value = -572159
$ws.Range("H97").Value = 1354.069
$ws.Range("I97").Value = 1335.25
$ws.Range("K97").Value = 1335.25
$ws.Range("M97").Value = -839.25
$ws.Range("H102").Value = 202835.8
$ws.Range("I102").Value = 288488.44
$ws.Range("J102").Value = 2979.6667
$ws.Range("K102").Value = 288488.44
$ws.Range("L102").Value = 2979.6667
$ws.Range("M102").Value = -286866.44
$ws.Range("N102").Value = -6223.6667
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 13409.13
$ws.Range("I86").Value = 6881.6
$ws.Range("J86").Value = 25433.525
$ws.Range("K86").Value = 6881.6
$ws.Range("L86").Value = 25433.525
$ws.Range("M86").Value = -5758.6
$ws.Range("N86").Value = -27679.525
$ws.Range("H89").Value = 13409.13
$ws.Range("I89").Value = 6881.6
$ws.Range("J89").Value = 25433.525
$ws.Range("K89").Value = 34408
$ws.Range("L89").Value = 127167.625
$ws.Range("M89").Value = -28792
$ws.Range("N89").Value = -138399.625
$ws.Range("H94").Value = 888.1786
$ws.Range("I94").Value = 373.75
$ws.Range("K94").Value = 373.75
$ws.Range("M94").Value = 77.25
$ws.Range("H99").Value = 4605.8335
$ws.Range("I99").Value = 1874.6666
$ws.Range("J99").Value = 7337
$ws.Range("K99").Value = 1874.6666
$ws.Range("L99").Value = 7337
$ws.Range("M99").Value = -376.6666
$ws.Range("N99").Value = -10333
$ws.Range("H134").Value = 2873.0667
$ws.Range("I134").Value = 2739.6
$ws.Range("K134").Value = 8218.799999999999
$ws.Range("M134").Value = -5683.799999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4417.433
$ws.Range("J31").Value = 5549.7856
$ws.Range("L31").Value = 5549.7856
$ws.Range("N31").Value = -6139.7856
$ws.Range("H34").Value = 4417.433
$ws.Range("J34").Value = 5549.7856
$ws.Range("L34").Value = 5549.7856
$ws.Range("N34").Value = -5953.7856
$ws.Range("H58").Value = 2560.4443
$ws.Range("I58").Value = 1982.4
$ws.Range("K58").Value = 1982.4
$ws.Range("M58").Value = -1779.4
$ws.Range("H103").Value = 18915.334
$ws.Range("I103").Value = 18915.334
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 18915.334
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -17743.334
$ws.Range("H105").Value = 70007
$ws.Range("J105").Value = 55005
$ws.Range("L105").Value = 55005
$ws.Range("N105").Value = -58499
$ws.Range("H134").Value = 1504
$ws.Range("I134").Value = 1006
$ws.Range("K134").Value = 3018
$ws.Range("M134").Value = -483
$ws.Range("H136").Value = 2560.4443
$ws.Range("I136").Value = 1982.4
$ws.Range("K136").Value = 5947.200000000001
$ws.Range("M136").Value = -3397.200000000001
$ws.Range("H141").Value = 254443.67
$ws.Range("J141").Value = 254443.67
$ws.Range("L141").Value = 254443.67
$ws.Range("N141").Value = -264803.67
$ws.Range("N103").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 11190.111
$ws.Range("I14").Value = 11190.111
$ws.Range("K14").Value = 33570.333
$ws.Range("M14").Value = -33397.333
$ws.Range("H122").Value = 734.7143
$ws.Range("J122").Value = 810.8
$ws.Range("L122").Value = 7297.2
$ws.Range("N122").Value = -12197.2
$ws.Range("H127").Value = 3987.25
$ws.Range("J127").Value = 3987.25
$ws.Range("L127").Value = 11961.75
$ws.Range("N127").Value = -21881.75
$ws.Range("H131").Value = 1566.3429
$ws.Range("J131").Value = 1692.25
$ws.Range("L131").Value = 5076.75
$ws.Range("N131").Value = -15156.75
$ws.Range("H134").Value = 3335.76
$ws.Range("I134").Value = 1114
$ws.Range("K134").Value = 3342
$ws.Range("M134").Value = 1728
$ws.Range("H136").Value = 2144.125
$ws.Range("I136").Value = 2144.125
$ws.Range("K136").Value = 6432.375
$ws.Range("M136").Value = -1332.375
$ws.Range("H137").Value = 5841.5
$ws.Range("I137").Value = 3849.5
$ws.Range("J137").Value = 6505.5
$ws.Range("K137").Value = 11548.5
$ws.Range("L137").Value = 19516.5
$ws.Range("M137").Value = -6448.5
$ws.Range("N137").Value = -29716.5
$ws.Range("H138").Value = 4450.6
$ws.Range("I138").Value = 4805
$ws.Range("K138").Value = 14415
$ws.Range("M138").Value = -9275
$ws.Range("H139").Value = 3534
$ws.Range("I139").Value = 3487.8
$ws.Range("K139").Value = 10463.4
$ws.Range("M139").Value = -5323.400000000001
$ws.Range("H141").Value = 9353.299999999999
$ws.Range("I141").Value = 7935.6
$ws.Range("J141").Value = 13606.4
$ws.Range("K141").Value = 23806.8
$ws.Range("L141").Value = 40819.2
$ws.Range("M141").Value = -18626.8
$ws.Range("N141").Value = -51179.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1725.6154
$ws.Range("I97").Value = 1251.1482
$ws.Range("J97").Value = 2793.1667
$ws.Range("K97").Value = 1251.1482
$ws.Range("L97").Value = 2793.1667
$ws.Range("M97").Value = -755.1482000000001
$ws.Range("N97").Value = -3785.1667

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2643.7188
$ws.Range("I82").Value = 1724.875
$ws.Range("J82").Value = 3562.5625
$ws.Range("K82").Value = 1724.875
$ws.Range("L82").Value = 3562.5625
$ws.Range("M82").Value = -1363.875
$ws.Range("N82").Value = -4284.5625
$ws.Range("H85").Value = 2643.7188
$ws.Range("I85").Value = 1724.875
$ws.Range("J85").Value = 3562.5625
$ws.Range("K85").Value = 1724.875
$ws.Range("L85").Value = 3562.5625
$ws.Range("M85").Value = -476.875
$ws.Range("N85").Value = -6058.5625
$ws.Range("H93").Value = 1843.4828
$ws.Range("I93").Value = 1097.7778
$ws.Range("J93").Value = 3063.7273
$ws.Range("K93").Value = 1097.7778
$ws.Range("L93").Value = 3063.7273
$ws.Range("M93").Value = 150.2221999999999
$ws.Range("N93").Value = -5559.7273
$ws.Range("H100").Value = 3429.3713
$ws.Range("I100").Value = 3061.5
$ws.Range("K100").Value = 3061.5
$ws.Range("M100").Value = -2520.5
$ws.Range("H136").Value = 8600.25
$ws.Range("I136").Value = 4701
$ws.Range("K136").Value = 14103
$ws.Range("M136").Value = -11553

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 29500
$ws.Range("I39").Value = 29500
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 29500
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -29087
$ws.Range("H132").Value = 3797
$ws.Range("I132").Value = 5376.4375
$ws.Range("J132").Value = 989.1111
$ws.Range("K132").Value = 16129.3125
$ws.Range("L132").Value = 2967.3333
$ws.Range("M132").Value = -13599.3125
$ws.Range("N132").Value = -8027.3333
$ws.Range("H136").Value = 3482.2917
$ws.Range("I136").Value = 1031.4667
$ws.Range("K136").Value = 3094.4001
$ws.Range("M136").Value = -544.4000999999998
$ws.Range("N39").ClearContents()

